# Adds Corequisites/Concurrent/Recommended columns (D:F), shifting the old
# "Terms Typically Offered" data into column G, and splits any prerequisite
# text that had embedded Concurrent/Corequisite/Recommended notes into the
# correct new column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Cells.Item(1, 4).Value = 'Corequisites'
$ws.Cells.Item(1, 5).Value = 'Concurrent'
$ws.Cells.Item(1, 6).Value = 'Recommended'
$ws.Cells.Item(1, 7).Value = 'Terms Typically Offered'

# --- A few Prerequisites (column C) cells had embedded "Recommended:" notes
#     removed now that there is a dedicated Recommended column, and two cells
#     had minor wording fixed ("one of the following:" -> "one of the"). ---
$prereqUpdates = @(
    @{ Row = 40; C = 'CSC/CPE 102 and CSC/CPE 103, or CSC/CPE 202 and CSC/CPE 203; and one of the STAT 301, STAT 312, STAT 321 or STAT 350.' }
    @{ Row = 42; C = 'CPE/CSC 202.' }
    @{ Row = 64; C = 'CSC 349 and one of the STAT 302, STAT 312, STAT 321 or STAT 350.' }
    @{ Row = 75; C = 'CSC 466 or CSC 480 or graduate standing.' }
    @{ Row = 79; C = 'CSC 349, and MATH 206 or MATH 244.' }
    @{ Row = 106; C = 'CSC 482 and graduate standing.' }
)
foreach ($u in $prereqUpdates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}

# --- Corequisites (D), Concurrent (E), Recommended (F) and Terms Typically
#     Offered (G, formerly column D) for every course row. ---
$rows = @(
    @{ Row = 2; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 3; D = 'NA'; E = 'CPE/CSC 101.'; F = 'NA'; G = 'TBD' }
    @{ Row = 4; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 5; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 6; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 7; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 8; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 9; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 10; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 11; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 12; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 13; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W' }
    @{ Row = 14; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 15; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 16; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 17; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 18; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 19; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 20; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 21; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 22; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 23; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, SP' }
    @{ Row = 24; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W' }
    @{ Row = 25; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W, SP' }
    @{ Row = 26; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 27; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 28; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 29; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 30; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 31; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 32; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 33; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 34; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 35; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, SP' }
    @{ Row = 36; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 37; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 38; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 39; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 40; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W' }
    @{ Row = 41; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W, SP' }
    @{ Row = 42; D = 'NA'; E = 'NA'; F = 'ART 376.'; G = 'SP ' }
    @{ Row = 43; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 44; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 45; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 46; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 47; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 48; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 49; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 50; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 51; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 52; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 53; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 54; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 55; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 56; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 57; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 58; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 59; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 60; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 61; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 62; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 63; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 64; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, SP' }
    @{ Row = 65; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 66; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 67; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 68; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 69; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 70; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 71; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 72; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 73; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, SP' }
    @{ Row = 74; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 75; D = 'NA'; E = 'NA'; F = 'CSC 349.'; G = 'F, W ' }
    @{ Row = 76; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 77; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 78; D = 'CSC 484.'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 79; D = 'NA'; E = 'NA'; F = 'CSC 466 or CSC 480.'; G = 'F ' }
    @{ Row = 80; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 81; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 82; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 83; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 84; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 85; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 86; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 87; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 88; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 89; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 90; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 91; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 92; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 93; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 94; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 95; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 96; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 97; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 98; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 99; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 100; D = 'NA'; E = 'NA'; F = 'NA'; G = 'SP' }
    @{ Row = 101; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 102; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W' }
    @{ Row = 103; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F' }
    @{ Row = 104; D = 'NA'; E = 'NA'; F = 'NA'; G = 'W' }
    @{ Row = 105; D = 'NA'; E = 'NA'; F = 'NA'; G = 'TBD' }
    @{ Row = 106; D = 'NA'; E = 'NA'; F = 'CSC 580.'; G = 'TBD ' }
    @{ Row = 107; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 108; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 109; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 110; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 111; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 112; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
    @{ Row = 113; D = 'NA'; E = 'NA'; F = 'NA'; G = 'F, W, SP' }
)

foreach ($item in $rows) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
    $ws.Cells.Item($item.Row, 5).Value = $item.E
    $ws.Cells.Item($item.Row, 6).Value = $item.F
    $ws.Cells.Item($item.Row, 7).Value = $item.G
}

